# Insert a new data row at row 361 (pushing the existing rows 361-483 down
# to 362-484) and populate the new row with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("361:361").Insert()

$ws.Range("A361").Value = 7
$ws.Range("B361").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C361").Value = 'Ñuble'
$ws.Range("D361").Value = 44900
$ws.Range("E361").Value = 16
$ws.Range("F361").Value = 100114014
$ws.Range("G361").Value = 'Betarraga'
$ws.Range("H361").Value = 'Sin especificar'
$ws.Range("I361").Value = 'Primera'
$ws.Range("J361").Value = 400
$ws.Range("K361").Value = 800
$ws.Range("L361").Value = 850
$ws.Range("M361").Value = 825
$ws.Range("N361").Value = '$/paquete 5 unidades'
$ws.Range("O361").Value = 'Provincia de Diguillín'
$ws.Range("P361").Value = 165
$ws.Range("Q361").Value = 5
$ws.Range("R361").Value = 'Hortaliza'
